$d = $word.ActiveDocument

function Get-ParaIndexByText($doc, $text) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text.TrimEnd([char]13) -eq $text) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 1) Title: "Committee <Number in Roman Numeral>" -> "Committee <Number in Roman Ordinals>"
# ---------------------------------------------------------------------------
[void]$d.Content.Find.Execute(
    "Committee <Number in Roman Numeral>", $true, $false, $false, $false, $false,
    $true, 1, $false, "Committee <Number in Roman Ordinals>", 2)

# ---------------------------------------------------------------------------
# 2) Insert "Committee Members" block right before the "Committee Assignment"
#    paragraph. NOTE: paragraph references are re-fetched by index/text after
#    every document mutation, since stale references can end up pointing at
#    the wrong paragraph once new paragraphs are spliced in.
# ---------------------------------------------------------------------------

# -- new paragraph: "Committee Members" (centered, bold heading, same style as
#    the "Committee Assignment" heading it precedes)
$caIndex = Get-ParaIndexByText $d "Committee Assignment"
[void]$d.Paragraphs.Item($caIndex).Range.InsertParagraphBefore()
$caIndex = Get-ParaIndexByText $d "Committee Assignment"
$membersHeading = $d.Paragraphs.Item($caIndex - 1)
$membersHeading.Alignment = 1
$membersHeading.Range.Text = "Committee Members"
$membersHeading.Range.Font.Name = "Roboto"
$membersHeading.Range.Font.Bold = -1

# -- new paragraph: the member list (left aligned, bold, with manual line
#    breaks between each member)
$caIndex = Get-ParaIndexByText $d "Committee Assignment"
[void]$d.Paragraphs.Item($caIndex).Range.InsertParagraphBefore()
$caIndex = Get-ParaIndexByText $d "Committee Assignment"
$membersList = $d.Paragraphs.Item($caIndex - 1)
$membersList.Alignment = 0
$memberText = "The members of this committee are:" + [char]11 + "<Member>," + [char]11 + "<Member>," + [char]11 + "and <Member>."
$membersList.Range.Text = $memberText
$membersList.Range.Font.Name = "Roboto"
$membersList.Range.Font.Bold = -1

# -- new paragraph: blank centered line
$caIndex = Get-ParaIndexByText $d "Committee Assignment"
[void]$d.Paragraphs.Item($caIndex).Range.InsertParagraphBefore()
$caIndex = Get-ParaIndexByText $d "Committee Assignment"
$blankLine = $d.Paragraphs.Item($caIndex - 1)
$blankLine.Alignment = 1

# ---------------------------------------------------------------------------
# 3) Insert "Expected Committee Report Format" block right before the
#    "Expected Completion Date" paragraph.
# ---------------------------------------------------------------------------

# -- new paragraph: "Expected Committee Report Format" (centered, bold heading)
$ecdIndex = Get-ParaIndexByText $d "Expected Completion Date"
[void]$d.Paragraphs.Item($ecdIndex).Range.InsertParagraphBefore()
$ecdIndex = Get-ParaIndexByText $d "Expected Completion Date"
$formatHeading = $d.Paragraphs.Item($ecdIndex - 1)
$formatHeading.Alignment = 1
$formatHeading.Range.Text = "Expected Committee Report Format"
$formatHeading.Range.Font.Name = "Roboto"
$formatHeading.Range.Font.Bold = -1

# -- new paragraph: description text (left aligned, bold)
$ecdIndex = Get-ParaIndexByText $d "Expected Completion Date"
[void]$d.Paragraphs.Item($ecdIndex).Range.InsertParagraphBefore()
$ecdIndex = Get-ParaIndexByText $d "Expected Completion Date"
$formatBody = $d.Paragraphs.Item($ecdIndex - 1)
$formatBody.Alignment = 0
$formatBody.Range.Text = "The expected format of the report of this committee is..."
$formatBody.Range.Font.Name = "Roboto"
$formatBody.Range.Font.Bold = -1

Write-Output "done"
